$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-24 00:43:58"

$wsZhCn.Range("H4").Value = "2016-08-24 00:43:53"
$wsZhCn.Range("K4").Value = "2016-08-24 00:44:15"

$wsDeDe.Range("K4").Value = "2016-08-24 00:44:22"
